$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览 (Exhibition)
$ws1.Range("F3").Value = 1460
$ws1.Range("F7").Value = 1274
$ws1.Range("F8").Value = 1628
$ws1.Range("F10").Value = 57
$ws1.Range("F11").Value = 2303
$ws1.Range("F12").Value = 470
$ws1.Range("F13").Value = 133
$ws1.Range("F16").Value = 97
$ws1.Range("F18").Value = 8306
$ws1.Range("F19").Value = 49
$ws1.Range("F20").Value = 6323
$ws1.Range("F21").Value = 10293
$ws1.Range("F24").Value = 189
$ws1.Range("F26").Value = 512
$ws1.Range("F28").Value = 157
$ws1.Range("F29").Value = 4409
$ws1.Range("F30").Value = 198
$ws1.Range("F31").Value = 404

$ws2 = $wb.Worksheets.Item(2)   # 演出 (Performance)
$ws2.Range("F8").Value = 1163
$ws2.Range("F20").Value = 12

$ws4 = $wb.Worksheets.Item(4)   # 全部类型 (All types)
$ws4.Range("F6").Value = 1460
$ws4.Range("F7").Value = 1460
$ws4.Range("F12").Value = 1274
$ws4.Range("F14").Value = 1628
$ws4.Range("F17").Value = 2303
$ws4.Range("F19").Value = 470
$ws4.Range("F20").Value = 133
$ws4.Range("F24").Value = 97
$ws4.Range("F26").Value = 8306
$ws4.Range("F27").Value = 49
$ws4.Range("F28").Value = 6323
$ws4.Range("F29").Value = 10293
$ws4.Range("F33").Value = 189
$ws4.Range("F36").Value = 512
$ws4.Range("F40").Value = 157
$ws4.Range("F41").Value = 4409
$ws4.Range("F43").Value = 198
$ws4.Range("F48").Value = 404
$ws4.Range("F49").Value = 12
